$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 279 ("搭乗券" / boarding pass row),
# shifting all subsequent rows down by one.
$ws.Rows.Item(279).Insert()

# New vocabulary entry: "パスポート" (passport) in both the Kanji/Kana columns,
# matching the heading-style formatting (Yu Gothic font, 18.75pt row height)
# used by the other single-column heading rows in this sheet.
$ws.Range("A279").Value = "パスポート"
$ws.Range("B279").Value = "パスポート"
$ws.Range("A279").Font.Name = "Yu Gothic"
$ws.Range("B279").Font.Name = "Yu Gothic"
$ws.Rows.Item(279).RowHeight = 18.75

# Update the current selection to match the author's cursor position after
# the edit.
$ws.Range("C279").Select()
